$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'319.04"
$ws.Range("E2").Value = "'3.67%"
$ws.Range("D3").Value = "'48.71"
$ws.Range("E3").Value = "'10.22%"
$ws.Range("D4").Value = "'5.314"
$ws.Range("E4").Value = "'4.34%"
$ws.Range("D5").Value = "'0.07934"
$ws.Range("E5").Value = "'-0.64%"
$ws.Range("D6").Value = "'4.584"
$ws.Range("E6").Value = "'3.44%"
$ws.Range("D7").Value = "'1.334"
$ws.Range("E7").Value = "'24.56%"
$ws.Range("D8").Value = "'1.623"
$ws.Range("E8").Value = "'1.50%"
$ws.Range("D9").Value = "'0.1252"
$ws.Range("E9").Value = "'-2.75%"
$ws.Range("D10").Value = "'0.1967"
$ws.Range("E10").Value = "'4.08%"
$ws.Range("D11").Value = "'0.09511"
$ws.Range("E11").Value = "'3.56%"
$ws.Range("D12").Value = "'0.04555"
$ws.Range("E12").Value = "'7.96%"
$ws.Range("D13").Value = "'0.1050"
$ws.Range("E13").Value = "'1.00%"
$ws.Range("D14").Value = "'0.001333"
$ws.Range("E14").Value = "'2.98%"
$ws.Range("D15").Value = "'0.04212"
$ws.Range("E15").Value = "'1.73%"
$ws.Range("D16").Value = "'0.005949"
$ws.Range("E16").Value = "'4.01%"
$ws.Range("D17").Value = "'3.345"
$ws.Range("E17").Value = "'0.13%"
$ws.Range("D18").Value = "'2.436"
$ws.Range("E18").Value = "'1.76%"
$ws.Range("D19").Value = "'0.3465"
$ws.Range("E19").Value = "'3.48%"
$ws.Range("D20").Value = "'8.135"
$ws.Range("E20").Value = "'1.11%"
$ws.Range("D21").Value = "'0.1403"
$ws.Range("E21").Value = "'3.05%"
$ws.Range("D22").Value = "'0.3076"
$ws.Range("E22").Value = "'10.29%"
$ws.Range("E23").Value = "'2.74%"
$ws.Range("D24").Value = "'0.004194"
$ws.Range("E24").Value = "'-2.93%"
$ws.Range("D25").Value = "'0.0001360"
$ws.Range("E25").Value = "'1.99%"
$ws.Range("D26").Value = "'0.0003564"
$ws.Range("E26").Value = "'-95.19%"
$ws.Range("D38").Value = "'0.02656"
$ws.Range("E38").Value = "'0.64%"
$ws.Range("D39").Value = "'0.05848"
$ws.Range("E39").Value = "'8.87%"
$ws.Range("D40").Value = "'0.01041"
$ws.Range("E40").Value = "'86.07%"
$ws.Range("D41").Value = "'0.008042"
$ws.Range("E41").Value = "'4.11%"
$ws.Range("D42").Value = "'0.1455"
$ws.Range("E42").Value = "'3.70%"
$ws.Range("D43").Value = "'0.007565"
$ws.Range("E43").Value = "'4.33%"
$ws.Range("D44").Value = "'0.007958"
$ws.Range("E44").Value = "'-5.27%"
$ws.Range("D45").Value = "'0.3196"
$ws.Range("E45").Value = "'4.20%"
$ws.Range("D46").Value = "'0.00007045"
$ws.Range("E46").Value = "'6.05%"
$ws.Range("D47").Value = "'0.00000000755"
$ws.Range("E47").Value = "'1.98%"
$ws.Range("E48").Value = "'9.36%"
$ws.Range("D49").Value = "'0.004028"
$ws.Range("E49").Value = "'2.00%"
$ws.Range("D50").Value = "'0.00002114"
$ws.Range("E50").Value = "'1.98%"
$ws.Range("D51").Value = "'0.0002014"
$ws.Range("E51").Value = "'1.98%"
